$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number + report date range) ---
$ws.Range("A8").Characters(21, 2).Text = "43"
$ws.Range("C9").Characters(27, 10).Text = "10/23/2023"
$ws.Range("C9").Characters(48, 10).Text = "10/29/2023"

# --- Cells that switch between numeric and placeholder-text ("0" / "***.*") ---
# Copy() from an untouched donor cell of the right style class first, so the
# destination lands on the SAME cellXf index Excel would reuse, then overwrite
# the value/text.
$ws.Range("C30").Copy($ws.Range("F14"))
$ws.Range("C30").Copy($ws.Range("D22"))
$ws.Range("E30").Copy($ws.Range("E22"))
$ws.Range("C30").Copy($ws.Range("G23"))
$ws.Range("E30").Copy($ws.Range("H23"))
$ws.Range("I30").Copy($ws.Range("D26"))
$ws.Range("K30").Copy($ws.Range("E26"))
$ws.Range("C30").Copy($ws.Range("D27"))
$ws.Range("E30").Copy($ws.Range("E27"))

# --- Final cell values ---
# Row 14
$ws.Range("F14").Value = "0"
$ws.Range("H14").Value = -100
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = -80
# Row 15
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 19
$ws.Range("K15").Value = -26.923076923076
$ws.Range("L15").Value = -5
$ws.Range("M15").Value = -13.636363636363
$ws.Range("N15").Value = -72.463768115942
# Row 16
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 11
$ws.Range("E16").Value = -18.181818181818
$ws.Range("F16").Value = 43
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = 72
$ws.Range("I16").Value = 296
$ws.Range("J16").Value = 326
$ws.Range("K16").Value = -9.20245398773
$ws.Range("L16").Value = -1.003344481605
$ws.Range("M16").Value = -18.681318681318
$ws.Range("N16").Value = -83.472920156337
# Row 17
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -15.384615384615
$ws.Range("F17").Value = 53
$ws.Range("G17").Value = 50
$ws.Range("H17").Value = 6
$ws.Range("I17").Value = 665
$ws.Range("J17").Value = 611
$ws.Range("K17").Value = 8.837970540098
$ws.Range("L17").Value = 33
$ws.Range("M17").Value = 87.323943661971
$ws.Range("N17").Value = -25.448430493273
# Row 18
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 33.333333333333
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -45
$ws.Range("I18").Value = 138
$ws.Range("J18").Value = 184
$ws.Range("K18").Value = -25
$ws.Range("L18").Value = -3.496503496503
$ws.Range("M18").Value = -21.142857142857
$ws.Range("N18").Value = -92.827442827442
# Row 19
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = 25
$ws.Range("F19").Value = 51
$ws.Range("G19").Value = 50
$ws.Range("H19").Value = 2
$ws.Range("I19").Value = 463
$ws.Range("J19").Value = 508
$ws.Range("K19").Value = -8.858267716535
$ws.Range("L19").Value = -10.271317829457
$ws.Range("M19").Value = 84.462151394422
$ws.Range("N19").Value = -43.467643467643
# Row 20
$ws.Range("C20").Value = 5
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 22
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = 46.666666666666
$ws.Range("I20").Value = 237
$ws.Range("J20").Value = 180
$ws.Range("K20").Value = 31.666666666666
$ws.Range("L20").Value = 85.15625
$ws.Range("M20").Value = 125.714285714286
$ws.Range("N20").Value = -68.4
# Row 21
$ws.Range("C21").Value = 45
$ws.Range("D21").Value = 44
$ws.Range("E21").Value = 2.272727272727
$ws.Range("F21").Value = 183
$ws.Range("G21").Value = 164
$ws.Range("H21").Value = 11.585365853658
$ws.Range("I21").Value = 1829
$ws.Range("J21").Value = 1847
$ws.Range("K21").Value = -0.974553329723
$ws.Range("L21").Value = 13.110698824984
$ws.Range("M21").Value = 42.556508183943
$ws.Range("N21").Value = -70.968253968254
# Row 22
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = "0"
$ws.Range("E22").Value = "***.*"
$ws.Range("F22").Value = 6
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 50
$ws.Range("I22").Value = 26
$ws.Range("K22").Value = 4
$ws.Range("L22").Value = -16.129032258064
$ws.Range("M22").Value = -3.703703703703
# Row 23
$ws.Range("G23").Value = "0"
$ws.Range("H23").Value = "***.*"
$ws.Range("L23").Value = -39.473684210526
# Row 24
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 33.333333333333
$ws.Range("F24").Value = 87
$ws.Range("G24").Value = 86
$ws.Range("H24").Value = 1.162790697674
$ws.Range("I24").Value = 882
$ws.Range("J24").Value = 1143
$ws.Range("K24").Value = -22.834645669291
$ws.Range("L24").Value = -9.259259259259
$ws.Range("M24").Value = 13.076923076923
# Row 25
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 26
$ws.Range("E25").Value = -42.307692307692
$ws.Range("F25").Value = 90
$ws.Range("G25").Value = 75
$ws.Range("H25").Value = 20
$ws.Range("I25").Value = 1039
$ws.Range("J25").Value = 894
$ws.Range("K25").Value = 16.219239373601
$ws.Range("L25").Value = 46.544428772919
$ws.Range("M25").Value = -4.853479853479
# Row 26
$ws.Range("D26").Value = 1
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 7
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = 75
$ws.Range("I26").Value = 50
$ws.Range("J26").Value = 44
$ws.Range("K26").Value = 13.636363636363
$ws.Range("L26").Value = 56.25
# Row 27
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = "0"
$ws.Range("E27").Value = "***.*"
$ws.Range("I27").Value = 74
$ws.Range("K27").Value = -1.333333333333
$ws.Range("L27").Value = 10.447761194029
# Row 28
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 50
$ws.Range("L28").Value = 27.5
$ws.Range("M28").Value = 2
$ws.Range("N28").Value = -66.883116883116
# Row 29
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 50
$ws.Range("L29").Value = 8.108108108108
$ws.Range("M29").Value = -6.976744186046
$ws.Range("N29").Value = -71.830985915493
